$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supplier / contact header block
$ws.Range("B3").Value = "宁波瑾秀制刷科技有限公司"
$ws.Range("B4").Value = "电话: +86-574-27889688 传真: +86-574-27889677"

# Line-item description, quantity and unit price
$ws.Range("C7").Value = "2175，值直径0.8mm尼龙针+仿猪鬃，出锋17mm，`n产品颜色做204C，染头颜色204c"
$ws.Range("D7").Value = 25569.33362927083
$ws.Range("D7").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("E7").Value = 5.72

# Color-card filename
$ws.Range("F14").Value = "204C.png"

# Packing / QC notes
$ws.Range("A19").Value = "1：表面不得有污渍"
$ws.Range("A20").Value = "2：毛丝切面整齐"
$ws.Range("A21").Value = "3：外箱单边不得超过60cm，重量不得超过20公斤，装箱数不得超过100，不能有尾箱"
$ws.Range("A22").Value = "4：尼龙丝针注意不要大量明显弯曲"
$ws.Range("A23").Value = "5：染头不能有坑(2个以上），（有坑一律次品），不得有漏染"
$ws.Range("A24").Value = "6：染头不得粘连"
$ws.Range("A25").Value = "7：每个外箱需要贴2个标签FBA标签一个，货代标签一个"
